# Update cryptocurrency price/volume figures per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.523.27"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.924.72"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.35"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4809"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4049"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08191"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.008"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.78"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.918.36"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.084"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.50"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06888"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.533.49"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.661"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.97"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.185"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.165.65"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.90"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.386"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.086"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.47"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09582"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.593"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.562"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.385"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06347"
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02277"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.189"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5943"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.71"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.010"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.889"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1843"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.471"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.276"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.41"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07468"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5545"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.973"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.59"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.428"
$ws.Range("E51").Value = "  +1.18%  "
